$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old two-row header (rows 1 & 2) -----------------------
# Deleting row 1 twice shifts the data block (old rows 3:13) up to 1:11.
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()

# --- Insert a brand-new single header row in its place -----------------
$ws.Rows.Item(1).Insert()

# --- Populate the new header row ----------------------------------------
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# --- Give the numeric/unit header cells (F1:K1) the same look as the ----
# --- rest of the sheet's header-style text (Arial 9) --------------------
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# --- Restore the expected selection highlight ----------------------------
$null = $ws.Range("A2:K2").Select()

Write-Output "ok"
